$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# C8: "<DELETE>Abwesenheit Sylvester" -> "<DELETE Abwesenheit Sylvester>"
$ws.Range("C8").Value = "<DELETE Abwesenheit Sylvester>"

# B2/C2: drop the highlighted style (revert to plain default formatting)
$ws.Range("B2:C2").ClearFormats()

# C2: "xpath=//android.widget.GridView" -> "//android.widget.GridView"
$ws.Range("C2").Value = "//android.widget.GridView"

# Update the active selection shown on the sheet
$ws.Range("F4").Select()
